# Update cryptos list (price/volume + two row swaps) per commit diff
# Applies cell-by-cell text writes, forcing text storage (leading apostrophe)
# and resetting style back to Normal so no quotePrefix/number-format styling
# side effect is introduced on numeric-looking values (e.g. "0.386", "248.65").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'36.997.99"
$ws.Cells.Item(2, 4).Style = 'Normal'
$ws.Cells.Item(2, 5).Value = "'  +2.36%  "
$ws.Cells.Item(2, 5).Style = 'Normal'
$ws.Cells.Item(3, 4).Value = "'2.036.87"
$ws.Cells.Item(3, 4).Style = 'Normal'
$ws.Cells.Item(3, 5).Value = "'  +1.20%  "
$ws.Cells.Item(3, 5).Style = 'Normal'
$ws.Cells.Item(4, 5).Value = "'  +0.05%  "
$ws.Cells.Item(4, 5).Style = 'Normal'
$ws.Cells.Item(5, 4).Value = "'248.65"
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = "'  -1.06%  "
$ws.Cells.Item(5, 5).Style = 'Normal'
$ws.Cells.Item(6, 4).Value = "'0.638"
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = "'  -0.74%  "
$ws.Cells.Item(6, 5).Style = 'Normal'
$ws.Cells.Item(7, 4).Value = "'63.21"
$ws.Cells.Item(7, 4).Style = 'Normal'
$ws.Cells.Item(7, 5).Value = "'  +1.05%  "
$ws.Cells.Item(7, 5).Style = 'Normal'
$ws.Cells.Item(8, 5).Value = "'  +0.02%  "
$ws.Cells.Item(8, 5).Style = 'Normal'
$ws.Cells.Item(9, 4).Value = "'0.386"
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = "'  +4.58%  "
$ws.Cells.Item(9, 5).Style = 'Normal'
$ws.Cells.Item(10, 4).Value = "'58.23"
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = "'  -1.54%  "
$ws.Cells.Item(10, 5).Style = 'Normal'
$ws.Cells.Item(11, 4).Value = "'0.0798"
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = "'  +6.79%  "
$ws.Cells.Item(11, 5).Style = 'Normal'
$ws.Cells.Item(12, 5).Value = "'  -0.66%  "
$ws.Cells.Item(12, 5).Style = 'Normal'
$ws.Cells.Item(13, 4).Value = "'0.897"
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).Value = "'  -3.33%  "
$ws.Cells.Item(13, 5).Style = 'Normal'
$ws.Cells.Item(14, 4).Value = "'22.74"
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).Value = "'  +16.78%  "
$ws.Cells.Item(14, 5).Style = 'Normal'
$ws.Cells.Item(15, 4).Value = "'14.37"
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = "'  -3.49%  "
$ws.Cells.Item(15, 5).Style = 'Normal'
$ws.Cells.Item(16, 4).Value = "'2.333.87"
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).Value = "'  +1.10%  "
$ws.Cells.Item(16, 5).Style = 'Normal'
$ws.Cells.Item(17, 4).Value = "'5.53"
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Cells.Item(17, 5).Value = "'  +2.02%  "
$ws.Cells.Item(17, 5).Style = 'Normal'
$ws.Cells.Item(18, 4).Value = "'2.040.90"
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).Value = "'  +1.71%  "
$ws.Cells.Item(18, 5).Style = 'Normal'
$ws.Cells.Item(19, 4).Value = "'36.938.15"
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = "'  +2.45%  "
$ws.Cells.Item(19, 5).Style = 'Normal'
$ws.Cells.Item(20, 4).Value = "'72.09"
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = "'  +0.15%  "
$ws.Cells.Item(20, 5).Style = 'Normal'
$ws.Cells.Item(21, 4).Value = "'0.0₃0878"
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = "'  +2.71%  "
$ws.Cells.Item(21, 5).Style = 'Normal'
$ws.Cells.Item(22, 4).Value = "'5.38"
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = "'  +2.25%  "
$ws.Cells.Item(22, 5).Style = 'Normal'
$ws.Cells.Item(23, 4).Value = "'236.11"
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = "'  +0.89%  "
$ws.Cells.Item(23, 5).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = "'  -0.08%  "
$ws.Cells.Item(24, 5).Style = 'Normal'
$ws.Cells.Item(25, 4).Value = "'2.52"
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = "'  -7.84%  "
$ws.Cells.Item(25, 5).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = "'  +1.83%  "
$ws.Cells.Item(26, 5).Style = 'Normal'
$ws.Cells.Item(27, 4).Value = "'9.72"
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = "'  +1.73%  "
$ws.Cells.Item(27, 5).Style = 'Normal'
$ws.Cells.Item(28, 4).Value = "'159.05"
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = "'  -3.32%  "
$ws.Cells.Item(28, 5).Style = 'Normal'
$ws.Cells.Item(29, 4).Value = "'20.19"
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = "'  +2.97%  "
$ws.Cells.Item(29, 5).Style = 'Normal'
$ws.Cells.Item(30, 4).Value = "'0.133"
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = "'  +24.74%  "
$ws.Cells.Item(30, 5).Style = 'Normal'
$ws.Cells.Item(31, 4).Value = "'0.121"
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Cells.Item(31, 5).Value = "'  +0.86%  "
$ws.Cells.Item(31, 5).Style = 'Normal'
$ws.Cells.Item(32, 4).Value = "'5.06"
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = "'  -1.94%  "
$ws.Cells.Item(32, 5).Style = 'Normal'
$ws.Cells.Item(33, 4).Value = "'1.17"
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = "'  -2.62%  "
$ws.Cells.Item(33, 5).Style = 'Normal'
$ws.Cells.Item(34, 4).Value = "'0.0618"
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = "'  +2.45%  "
$ws.Cells.Item(34, 5).Style = 'Normal'
$ws.Cells.Item(35, 4).Value = "'4.50"
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).Value = "'  +0.26%  "
$ws.Cells.Item(35, 5).Style = 'Normal'
$ws.Cells.Item(36, 4).Value = "'2.37"
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = "'  -4.25%  "
$ws.Cells.Item(36, 5).Style = 'Normal'
$ws.Cells.Item(37, 4).Value = "'6.31"
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = "'  +10.22%  "
$ws.Cells.Item(37, 5).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = "'  +0.05%  "
$ws.Cells.Item(38, 5).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = "'  +1.23%  "
$ws.Cells.Item(39, 5).Style = 'Normal'
$ws.Cells.Item(40, 4).Value = "'3.06"
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = "'  +29.51%  "
$ws.Cells.Item(40, 5).Style = 'Normal'
$ws.Cells.Item(41, 2).Value = "'Cronos"
$ws.Cells.Item(41, 2).Style = 'Normal'
$ws.Cells.Item(41, 3).Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(41, 3).Style = 'Normal'
$ws.Cells.Item(41, 4).Value = "'0.0989"
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).Value = "'  -9.04%  "
$ws.Cells.Item(41, 5).Style = 'Normal'
$ws.Cells.Item(42, 2).Value = "'TrustWalletToken"
$ws.Cells.Item(42, 2).Style = 'Normal'
$ws.Cells.Item(42, 3).Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(42, 3).Style = 'Normal'
$ws.Cells.Item(42, 4).Value = "'1.24"
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = "'  +1.61%  "
$ws.Cells.Item(42, 5).Style = 'Normal'
$ws.Cells.Item(43, 4).Value = "'2.98"
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = "'  +3.40%  "
$ws.Cells.Item(43, 5).Style = 'Normal'
$ws.Cells.Item(44, 2).Value = "'InjectiveProtocol"
$ws.Cells.Item(44, 2).Style = 'Normal'
$ws.Cells.Item(44, 3).Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(44, 3).Style = 'Normal'
$ws.Cells.Item(44, 4).Value = "'17.05"
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = "'  +1.78%  "
$ws.Cells.Item(44, 5).Style = 'Normal'
$ws.Cells.Item(45, 2).Value = "'ARBITRUM"
$ws.Cells.Item(45, 2).Style = 'Normal'
$ws.Cells.Item(45, 3).Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(45, 3).Style = 'Normal'
$ws.Cells.Item(45, 4).Value = "'1.13"
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = "'  +1.05%  "
$ws.Cells.Item(45, 5).Style = 'Normal'
$ws.Cells.Item(46, 4).Value = "'0.0214"
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = "'  -0.52%  "
$ws.Cells.Item(46, 5).Style = 'Normal'
$ws.Cells.Item(47, 4).Value = "'93.42"
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = "'  -0.63%  "
$ws.Cells.Item(47, 5).Style = 'Normal'
$ws.Cells.Item(48, 4).Value = "'7.69"
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = "'  -1.61%  "
$ws.Cells.Item(48, 5).Style = 'Normal'
$ws.Cells.Item(49, 4).Value = "'1.365.23"
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = "'  -4.15%  "
$ws.Cells.Item(49, 5).Style = 'Normal'
$ws.Cells.Item(50, 4).Value = "'2.92"
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = "'  +0.63%  "
$ws.Cells.Item(50, 5).Style = 'Normal'
$ws.Cells.Item(51, 4).Value = "'2.221.96"
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).Value = "'  +1.13%  "
$ws.Cells.Item(51, 5).Style = 'Normal'
